# "se agrega pdf del tp" -- adds the repo URL as its own paragraph right
# after the existing (empty) paragraph that follows "Repositorio;".

$d = $word.ActiveDocument

# Locate the "Repositorio;" paragraph, then the empty paragraph right
# after it -- that's where the new paragraph needs to be inserted.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Repositorio;`r") {
        $anchorIndex = $i
        break
    }
}

$emptyPara = $d.Paragraphs.Item($anchorIndex + 1)

# Insert a brand-new paragraph right after it and fill it with the URL;
# the new paragraph inherits the surrounding paragraph/run formatting
# (sz 24 / szCs 24, left indent 426 twips, exact 184-twip line spacing).
$emptyPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 2)
$newPara.Range.Text = "https://github.com/GabrielOkArg/TP_Composite_anegamiento.git"
